$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.963.13"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.353.81"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "240.24"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.669"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.73%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "74.10"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("E8").Value = "  -0.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.601"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("E10").Value = "  -0.06%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "60.35"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.49%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "33.28"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.40%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.26"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "2.704.78"
$ws.Range("E15").Value = "  -0.11%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "16.23"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.907"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "2.357.67"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "43.927.98"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  +0.13%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "78.19"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.85%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "253.27"
$c.Style = "Normal"
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.81"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -0.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.50"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "22.28"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -1.47%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0749"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  -3.69%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.49"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +14.98%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "64.77"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +13.55%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "19.27"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "9.17"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  -2.62%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.200"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -6.28%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -1.82%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.90%  "
$ws.Range("E50").Value = "  -2.15%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "98.65"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.96%  "
